$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-03-08"

# Update the header label cell (I1) that states the "through" date
$ws.Range("I1").Value = "2022 (through 03-08)"

# Update the data values that changed for March (row 4) and Total (row 14)
$ws.Range("I4").Value = 38
$ws.Range("I14").Value = 339
